# Generate Report for Handback
# This script reflects a localization handback report generation:
#  - Status cells move from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The zh-cn sheet gets its handback target/file/datetime columns (I, J, K) populated
#  - The de-de sheet gets its handback target/file/datetime columns (I, J, K) populated
#    with a later datetime than zh-cn
#  - A couple of columns are widened to better display the new hyperlink content

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Update status text everywhere it is used (Overview!E2:F3, and the
#    "Status" column (C) of the per-language report sheets).
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (I), "Latest Handback File" (J)
#    and "Latest Handback DateTime" (K) for both rows, and add hyperlinks on
#    the newly populated Target File cells (mirroring the Source File Name
#    hyperlinks already present in column A).
# ---------------------------------------------------------------------------
$zhHandbackDateTime = "2016-10-25 03:07:24"

$wsZhCn.Range("J2").Value = "5f266658-4f52-4a5f-9ade-256926e29f04.339593c65be1ddc4521d83985eb81b8f732de9cd.zh-cn.xlf"
$wsZhCn.Range("K2").Value = $zhHandbackDateTime

$wsZhCn.Range("J3").Value = "6d4abbe6-e83c-4be9-9a23-aa84760ac57b.b74c1aa282d0e352b610908cb4751e05180ac5cd.zh-cn.xlf"
$wsZhCn.Range("K3").Value = $zhHandbackDateTime

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e025622058dd88a41bf5ae265f9e63e081f58d2f/e2e/5f266658-4f52-4a5f-9ade-256926e29f04.md", "", "", "5f266658-4f52-4a5f-9ade-256926e29f04.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e025622058dd88a41bf5ae265f9e63e081f58d2f/e2e/6d4abbe6-e83c-4be9-9a23-aa84760ac57b.md", "", "", "6d4abbe6-e83c-4be9-9a23-aa84760ac57b.md")

# ---------------------------------------------------------------------------
# 3. de-de sheet: same as above but with its own (later) handback datetime.
# ---------------------------------------------------------------------------
$deHandbackDateTime = "2016-10-25 03:07:41"

$wsDeDe.Range("J2").Value = "5f266658-4f52-4a5f-9ade-256926e29f04.339593c65be1ddc4521d83985eb81b8f732de9cd.de-de.xlf"
$wsDeDe.Range("K2").Value = $deHandbackDateTime

$wsDeDe.Range("J3").Value = "6d4abbe6-e83c-4be9-9a23-aa84760ac57b.b74c1aa282d0e352b610908cb4751e05180ac5cd.de-de.xlf"
$wsDeDe.Range("K3").Value = $deHandbackDateTime

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e025622058dd88a41bf5ae265f9e63e081f58d2f/e2e/5f266658-4f52-4a5f-9ade-256926e29f04.md", "", "", "5f266658-4f52-4a5f-9ade-256926e29f04.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e025622058dd88a41bf5ae265f9e63e081f58d2f/e2e/6d4abbe6-e83c-4be9-9a23-aa84760ac57b.md", "", "", "6d4abbe6-e83c-4be9-9a23-aa84760ac57b.md")

# ---------------------------------------------------------------------------
# 4. Widen columns that now show the longer hyperlink / file-name text.
#    (ColumnWidth is expressed in characters; the values below are chosen so
#    the saved column width lands as close as possible to the target width.)
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth  = 29.166666666666668   # E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth  = 29.166666666666668   # F (de-de)

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.166666666666668       # C  Status
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.166666666666664       # I  Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664       # J  Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.166666666666668       # C  Status
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.166666666666664       # I  Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664       # J  Latest Handback File
